$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") "27.625.72"
Set-TextValue $ws.Range("E2") "  -1.58%  "
Set-TextValue $ws.Range("D3") "1.667.98"
Set-TextValue $ws.Range("E3") "  -3.20%  "
Set-TextValue $ws.Range("E4") "  +0.02%  "
Set-TextValue $ws.Range("D5") "215.23"
Set-TextValue $ws.Range("E5") "  -1.79%  "
Set-TextValue $ws.Range("D6") "0.509"
Set-TextValue $ws.Range("E6") "  -2.30%  "
Set-TextValue $ws.Range("E7") "  +0.00%  "
Set-TextValue $ws.Range("D8") "24.05"
Set-TextValue $ws.Range("E8") "  -2.31%  "
Set-TextValue $ws.Range("E9") "  -0.55%  "
Set-TextValue $ws.Range("E10") "  -1.70%  "
Set-TextValue $ws.Range("E11") "  -2.11%  "
Set-TextValue $ws.Range("D12") "1.902.60"
Set-TextValue $ws.Range("E12") "  -3.27%  "
Set-TextValue $ws.Range("D13") "1.667.43"
Set-TextValue $ws.Range("E13") "  -3.21%  "
Set-TextValue $ws.Range("D14") "4.14"
Set-TextValue $ws.Range("E14") "  -3.20%  "
Set-TextValue $ws.Range("D15") "0.559"
Set-TextValue $ws.Range("E15") "  -0.49%  "
Set-TextValue $ws.Range("D16") "66.63"
Set-TextValue $ws.Range("E16") "  -1.39%  "
Set-TextValue $ws.Range("D17") "27.605.23"
Set-TextValue $ws.Range("E17") "  -1.50%  "
Set-TextValue $ws.Range("D18") "243.62"
Set-TextValue $ws.Range("E18") "  +0.34%  "
Set-TextValue $ws.Range("D19") "0.0₃0732"
Set-TextValue $ws.Range("E19") "  -3.05%  "
Set-TextValue $ws.Range("D20") "7.67"
Set-TextValue $ws.Range("E20") "  -4.35%  "
Set-TextValue $ws.Range("E21") "  -0.04%  "
Set-TextValue $ws.Range("D22") "4.50"
Set-TextValue $ws.Range("E22") "  -2.97%  "
Set-TextValue $ws.Range("D23") "9.33"
Set-TextValue $ws.Range("E23") "  -3.97%  "
Set-TextValue $ws.Range("D24") "2.04"
Set-TextValue $ws.Range("E24") "  -4.40%  "
Set-TextValue $ws.Range("D25") "147.20"
Set-TextValue $ws.Range("E25") "  -1.11%  "
Set-TextValue $ws.Range("D26") "7.21"
Set-TextValue $ws.Range("E26") "  -3.71%  "
Set-TextValue $ws.Range("D27") "16.49"
Set-TextValue $ws.Range("E27") "  -1.39%  "
Set-TextValue $ws.Range("E28") "  +0.15%  "
Set-TextValue $ws.Range("E29") "  -2.41%  "
Set-TextValue $ws.Range("E30") "  +3.20%  "
Set-TextValue $ws.Range("E32") "  -2.41%  "
Set-TextValue $ws.Range("D33") "1.472.76"
Set-TextValue $ws.Range("E33") "  -1.47%  "
Set-TextValue $ws.Range("D34") "3.12"
Set-TextValue $ws.Range("E34") "  -4.85%  "
Set-TextValue $ws.Range("E35") "  -5.01%  "
Set-TextValue $ws.Range("E36") "  -1.08%  "
Set-TextValue $ws.Range("D37") "0.931"
Set-TextValue $ws.Range("E37") "  -2.67%  "
Set-TextValue $ws.Range("D38") "0.577"
Set-TextValue $ws.Range("E38") "  -5.10%  "
Set-TextValue $ws.Range("E39") "  -1.57%  "
Set-TextValue $ws.Range("D40") "69.65"
Set-TextValue $ws.Range("E40") "  -1.51%  "
Set-TextValue $ws.Range("E41") "  -4.55%  "
Set-TextValue $ws.Range("E42") "  -0.02%  "
Set-TextValue $ws.Range("B43") "FraxShare"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D43") "5.41"
Set-TextValue $ws.Range("E43") "  -7.45%  "
Set-TextValue $ws.Range("B44") "MXToken"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D44") "2.22"
Set-TextValue $ws.Range("E44") "  -2.91%  "
Set-TextValue $ws.Range("D45") "1.811.44"
Set-TextValue $ws.Range("E45") "  -3.14%  "
Set-TextValue $ws.Range("D46") "0.788"
Set-TextValue $ws.Range("E46") "  -1.74%  "
Set-TextValue $ws.Range("D47") "1.73"
Set-TextValue $ws.Range("E47") "  -1.95%  "
Set-TextValue $ws.Range("D48") "89.28"
Set-TextValue $ws.Range("E48") "  -1.91%  "
Set-TextValue $ws.Range("E49") "  -4.19%  "
Set-TextValue $ws.Range("E50") "  -2.26%  "
Set-TextValue $ws.Range("D51") "7.90"
Set-TextValue $ws.Range("E51") "  -4.49%  "
